# Calibration Legs and Others Update!
# Update the raw calibration input values (columns B-G) on Sheet1 for rows 4-9.
# Formula cells (H-M) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (LEG 1)
$ws.Range("B4").Value = 1450
$ws.Range("C4").Value = 1550

# Row 5 (LEG 2)
$ws.Range("B5").Value = 1400
$ws.Range("C5").Value = 1600
$ws.Range("D5").Value = 1600
$ws.Range("E5").Value = 1000
$ws.Range("F5").Value = 2070
$ws.Range("G5").Value = 2000

# Row 6 (LEG 3)
$ws.Range("B6").Value = 1400
$ws.Range("C6").Value = 1400
$ws.Range("D6").Value = 1600
$ws.Range("F6").Value = 1850
$ws.Range("G6").Value = 2050

# Row 7 (LEG 4)
$ws.Range("D7").Value = 1380
$ws.Range("F7").Value = 940

# Row 8 (LEG 5)
$ws.Range("B8").Value = 1400
$ws.Range("D8").Value = 1200
$ws.Range("E8").Value = 1850
$ws.Range("F8").Value = 1040

# Row 9 (LEG 6)
$ws.Range("B9").Value = 1300
$ws.Range("E9").Value = 1800
$ws.Range("G9").Value = 850

# Recalculate all dependent formula cells (H:M)
$excel.Calculate()

# Restore the active selection to F12 as recorded after the edit session
$ws.Range("F12").Select()
